$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A97").Copy()
$ws.Range("A98").PasteSpecial(-4122)

$ws.Range("A98").Value = 46003
$ws.Range("B98").Value = "22,1547"
$ws.Range("C98").Value = "15,9245"
$ws.Range("D98").Value = "15,7674"
$ws.Range("E98").Value = "15,7674"
